# Add team record (Wins/Losses/Ties) columns to the 2008 STL roster sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row (row 1): AD1="Wins", AE1="Losses", AF1="Ties"
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Match the bold/centered/bordered header style used by the other header cells (copy from AC1).
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Fill in the team record for every data row (2-45): 86 wins, 76 losses, 0 ties.
for ($r = 2; $r -le 45; $r++) {
    $ws.Cells.Item($r, 30).Value = 86
    $ws.Cells.Item($r, 31).Value = 76
    $ws.Cells.Item($r, 32).Value = 0
}

Write-Host "Team record columns added"
